$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Tooltip"), shifting existing
# columns E:L to F:M. Excel copies formatting from the column to the
# left (D) into the newly inserted column, matching the target layout.
$dColumnWidth = $ws.Columns("D:D").ColumnWidth
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = $dColumnWidth

# Header for the new "Alias Export" column.
$ws.Range("E1").Value = "Alias Export"

# Machine-readable alias identifiers for each field row.
$ws.Range("E2").Value = "REFERENCED_REPORTS"
$ws.Range("E3").Value = "NUCLEAR_ENERGY_ACTIVITIES_SEC_4_26"
$ws.Range("E4").Value = "NUCLEAR_ENERGY_ACTIVITIES_SEC_4_27"
$ws.Range("E5").Value = "NUCLEAR_ENERGY_ACTIVITIES_SEC_4_28"
$ws.Range("E6").Value = "FOSSIL_GAS_ACTIVITIES_SEC_4_29"
$ws.Range("E7").Value = "FOSSIL_GAS_ACTIVITIES_SEC_4_30"
$ws.Range("E8").Value = "FOSSIL_GAS_ACTIVITIES_SEC_4_31"
$ws.Range("E9").Value = "REV_ALIGNED_DENOMINATOR"
$ws.Range("E10").Value = "CAPEX_ALIGNED_DENOMINATOR"
$ws.Range("E11").Value = "REV_ALIGNED_NUMERATOR"
$ws.Range("E12").Value = "CAPEX_ALIGNED_NUMERATOR"
$ws.Range("E13").Value = "REV_NON_ALIGNED"
$ws.Range("E14").Value = "CAPEX_NON_ALIGNED"
$ws.Range("E15").Value = "REV_NON_ELIGIBLE"
$ws.Range("E16").Value = "CAPEX_NON_ELIGIBLE"
